# comandos de branches agregados a comandos Git
#
# 1) The "git init" paragraph gets a w:proofErr spellStart/spellEnd pair
#    wrapped around the "git" run (in addition to the existing gramStart/
#    gramEnd pair that was already there).
# 2) The trailing "mas.." paragraph is replaced with a set of new
#    paragraphs documenting git branch commands, and a new empty
#    paragraph is appended at the very end of the document.

$d = $word.ActiveDocument
$wNs = 'http://schemas.openxmlformats.org/wordprocessingml/2006/main'

# --- Change 1: wrap "git" (in "git init") with spellStart/spellEnd too ---

$gitInitPara = $d.Paragraphs.Item(2)
$gitInitXml = '<w:p xmlns:w="' + $wNs + '">' +
    '<w:pPr><w:ind w:left="707"/></w:pPr>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>git</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>init</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '</w:p>'
$null = $gitInitPara.Range.InsertXML($gitInitXml)

# --- Change 2: replace "mas.." paragraph with the new Branches content ---

$lastPara = $d.Paragraphs.Last
$branchesXml =
    '<w:p xmlns:w="' + $wNs + '">' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:t>Branches</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
    '</w:p>' +
    '<w:p xmlns:w="' + $wNs + '">' +
        '<w:r><w:t xml:space="preserve">Creación de </w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:t>Branches</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:t xml:space="preserve"> (ramas)</w:t></w:r>' +
    '</w:p>' +
    '<w:p xmlns:w="' + $wNs + '">' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:t>Branch</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:t xml:space="preserve"> (Rama): Es una copia del proyecto, bajo el control de versiones, de forma que los cambios realizados en esta rama no afecten al resto del proyecto y viceversa</w:t></w:r>' +
    '</w:p>' +
    '<w:p xmlns:w="' + $wNs + '">' +
        '<w:pPr><w:ind w:left="707"/></w:pPr>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:proofErr w:type="gramStart"/>' +
        '<w:r><w:t>git</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '<w:proofErr w:type="gramEnd"/>' +
        '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:t>branch</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:t>nombre_branche</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
    '</w:p>' +
    '<w:p xmlns:w="' + $wNs + '">' +
        '<w:r><w:t xml:space="preserve">Navegación sobre los </w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:t>branches</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
    '</w:p>' +
    '<w:p xmlns:w="' + $wNs + '">' +
        '<w:r><w:tab/></w:r>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:t>git</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:t>checkout</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:t xml:space="preserve"> nombre_branche</w:t></w:r>' +
        '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
        '<w:bookmarkEnd w:id="0"/>' +
    '</w:p>' +
    '<w:p xmlns:w="' + $wNs + '"/>'

$null = $lastPara.Range.InsertXML($branchesXml)
